# Auto-generated edit script: replicate scheduled-runner market-data refresh
# across the Leve profit tables (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1632.5
$ws.Range("I33").Value = 1132.2727
$ws.Range("J33").Value = 3466.6667
$ws.Range("K33").Value = 1132.2727
$ws.Range("L33").Value = 3466.6667
$ws.Range("M33").Value = -903.2727
$ws.Range("N33").Value = -3924.6667
$ws.Range("H41").Value = 136
$ws.Range("I41").Value = 175
$ws.Range("J41").Value = 126.25
$ws.Range("K41").Value = 175
$ws.Range("L41").Value = 126.25
$ws.Range("M41").Value = 265
$ws.Range("N41").Value = -1006.25
$ws.Range("H74").Value = 4893.5
$ws.Range("J74").Value = 5887
$ws.Range("L74").Value = 5887
$ws.Range("N74").Value = -7759
$ws.Range("H77").Value = 4893.5
$ws.Range("J77").Value = 5887
$ws.Range("L77").Value = 29435
$ws.Range("N77").Value = -38795
$ws.Range("H116").Value = 3081.889
$ws.Range("I116").Value = 2809.0908
$ws.Range("J116").Value = 3510.5715
$ws.Range("K116").Value = 2809.0908
$ws.Range("L116").Value = 3510.5715
$ws.Range("M116").Value = 632.9092000000001
$ws.Range("N116").Value = -10394.5715
$ws.Range("H127").Value = 1294.6
$ws.Range("I127").Value = 686.75
$ws.Range("J127").Value = 1515.6364
$ws.Range("K127").Value = 2060.25
$ws.Range("L127").Value = 4546.9092
$ws.Range("M127").Value = 2899.75
$ws.Range("N127").Value = -14466.9092
$ws.Range("H138").Value = 2538.807
$ws.Range("I138").Value = 2892.6428
$ws.Range("J138").Value = 2471.8647
$ws.Range("K138").Value = 8677.928400000001
$ws.Range("L138").Value = 7415.5941
$ws.Range("M138").Value = -3537.928400000001
$ws.Range("N138").Value = -17695.5941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 368413.66
$ws.Range("I32").Value = 402295.1
$ws.Range("J32").Value = 25364.25
$ws.Range("K32").Value = 402295.1
$ws.Range("L32").Value = 25364.25
$ws.Range("M32").Value = -402008.1
$ws.Range("N32").Value = -25938.25
$ws.Range("H74").Value = 2578.4
$ws.Range("I74").Value = 2388.4546
$ws.Range("J74").Value = 3971.3333
$ws.Range("K74").Value = 2388.4546
$ws.Range("L74").Value = 3971.3333
$ws.Range("M74").Value = -1514.4546
$ws.Range("N74").Value = -5719.3333
$ws.Range("H77").Value = 2578.4
$ws.Range("I77").Value = 2388.4546
$ws.Range("J77").Value = 3971.3333
$ws.Range("K77").Value = 11942.273
$ws.Range("L77").Value = 19856.6665
$ws.Range("M77").Value = -7574.273000000001
$ws.Range("N77").Value = -28592.6665
$ws.Range("H132").Value = 5517.0713
$ws.Range("I132").Value = 5058.364
$ws.Range("J132").Value = 5813.8823
$ws.Range("K132").Value = 15175.092
$ws.Range("L132").Value = 17441.6469
$ws.Range("M132").Value = -12645.092
$ws.Range("N132").Value = -22501.6469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 125002710
$ws.Range("I86").Value = 125002710
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 125002710
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -125001587
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 125002710
$ws.Range("I89").Value = 125002710
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 625013550
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -625007934
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 988.7143
$ws.Range("I94").Value = 498.58334
$ws.Range("J94").Value = 1642.2222
$ws.Range("K94").Value = 498.58334
$ws.Range("L94").Value = 1642.2222
$ws.Range("M94").Value = -47.58334000000002
$ws.Range("N94").Value = -2544.2222
$ws.Range("H134").Value = 2545.1428
$ws.Range("I134").Value = 2319.15
$ws.Range("K134").Value = 6957.450000000001
$ws.Range("M134").Value = -4422.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 895.6667
$ws.Range("I16").Value = 908.8570999999999
$ws.Range("K16").Value = 908.8570999999999
$ws.Range("M16").Value = -621.8570999999999
$ws.Range("H58").Value = 1722.8096
$ws.Range("I58").Value = 1868.5
$ws.Range("J58").Value = 1664.5333
$ws.Range("K58").Value = 1868.5
$ws.Range("L58").Value = 1664.5333
$ws.Range("M58").Value = -1665.5
$ws.Range("N58").Value = -2070.5333
$ws.Range("H62").Value = 4480
$ws.Range("I62").Value = 4422.222
$ws.Range("K62").Value = 4422.222
$ws.Range("M62").Value = -3798.222
$ws.Range("H65").Value = 4480
$ws.Range("I65").Value = 4422.222
$ws.Range("K65").Value = 22111.11
$ws.Range("M65").Value = -18991.11
$ws.Range("H99").Value = 1923.8096
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 1945
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1945
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -4941
$ws.Range("H105").Value = 1995.4286
$ws.Range("I105").Value = 1993.6
$ws.Range("K105").Value = 1993.6
$ws.Range("M105").Value = -246.5999999999999
$ws.Range("H112").Value = 37500
$ws.Range("J112").Value = 37500
$ws.Range("L112").Value = 37500
$ws.Range("N112").Value = -40454
$ws.Range("H113").Value = 895.6667
$ws.Range("I113").Value = 908.8570999999999
$ws.Range("K113").Value = 908.8570999999999
$ws.Range("M113").Value = 1261.1429
$ws.Range("H126").Value = 1923.8096
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1945
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5835
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10775
$ws.Range("H136").Value = 1722.8096
$ws.Range("I136").Value = 1868.5
$ws.Range("J136").Value = 1664.5333
$ws.Range("K136").Value = 5605.5
$ws.Range("L136").Value = 4993.5999
$ws.Range("M136").Value = -3055.5
$ws.Range("N136").Value = -10093.5999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1031.7693
$ws.Range("I108").Value = 442.3
$ws.Range("J108").Value = 2996.6667
$ws.Range("K108").Value = 1326.9
$ws.Range("L108").Value = 8990.000100000001
$ws.Range("M108").Value = 1553.1
$ws.Range("N108").Value = -14750.0001
$ws.Range("H113").Value = 989.2692
$ws.Range("I113").Value = 552
$ws.Range("J113").Value = 1093.381
$ws.Range("K113").Value = 1656
$ws.Range("L113").Value = 3280.143
$ws.Range("M113").Value = 514
$ws.Range("N113").Value = -7620.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 42639.5
$ws.Range("J141").Value = 42639.5
$ws.Range("L141").Value = 42639.5
$ws.Range("N141").Value = -52999.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 62338
$ws.Range("J11").Value = 62338
$ws.Range("L11").Value = 62338
$ws.Range("N11").Value = -62618
$ws.Range("H132").Value = 2611
$ws.Range("I132").Value = 1673.3
$ws.Range("J132").Value = 4053.6155
$ws.Range("K132").Value = 5019.9
$ws.Range("L132").Value = 12160.8465
$ws.Range("M132").Value = -2489.9
$ws.Range("N132").Value = -17220.8465
$ws.Range("H136").Value = 12822637
$ws.Range("I136").Value = 1862.375
$ws.Range("J136").Value = 33335876
$ws.Range("K136").Value = 5587.125
$ws.Range("L136").Value = 100007628
$ws.Range("M136").Value = -3037.125
$ws.Range("N136").Value = -100012728

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 102033.336
$ws.Range("I62").Value = 2100
$ws.Range("K62").Value = 2100
$ws.Range("M62").Value = -1476
$ws.Range("H65").Value = 102033.336
$ws.Range("I65").Value = 2100
$ws.Range("K65").Value = 10500
$ws.Range("M65").Value = -7380
$ws.Range("H132").Value = 4388307.5
$ws.Range("I132").Value = 2471.75
$ws.Range("J132").Value = 9261458
$ws.Range("K132").Value = 7415.25
$ws.Range("L132").Value = 27784374
$ws.Range("M132").Value = -4885.25
$ws.Range("N132").Value = -27789434
$ws.Range("H135").Value = 58750.4
$ws.Range("J135").Value = 58750.4
$ws.Range("L135").Value = 58750.4
$ws.Range("N135").Value = -68890.39999999999
$ws.Range("H136").Value = 2230.7368
$ws.Range("I136").Value = 2027.2084
$ws.Range("J136").Value = 2579.6428
$ws.Range("K136").Value = 6081.6252
$ws.Range("L136").Value = 7738.928400000001
$ws.Range("M136").Value = -3531.6252
$ws.Range("N136").Value = -12838.9284
$ws.Range("H140").Value = 36582
$ws.Range("J140").Value = 36582
$ws.Range("L140").Value = 36582
$ws.Range("N140").Value = -46942
$ws.Range("H141").Value = 44000
$ws.Range("J141").Value = 44000
$ws.Range("L141").Value = 44000
$ws.Range("N141").Value = -54360
